# Insert a new row at 181, shifting existing rows 181-186 down to 182-187,
# then populate the new row 181 with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 181; this shifts rows 181..186 to 182..187
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with data
$ws.Range("A181").Value = 10
$ws.Range("B181").Value = "Vega Modelo de Temuco"
$ws.Range("C181").Value = "La Araucanía"
$ws.Range("D181").Value = "2023-08-09"
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = "Fruta"
$ws.Range("G181").Value = 100107
$ws.Range("H181").Value = "Otros"
$ws.Range("I181").Value = 100107002
$ws.Range("J181").Value = "Chirimoya"
$ws.Range("K181").Value = "Cultivar IV Región"
$ws.Range("L181").Value = "Primera"
$ws.Range("M181").Value = 70
$ws.Range("N181").Value = 3500
$ws.Range("O181").Value = 4000
$ws.Range("P181").Value = 3643
$ws.Range("Q181").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R181").Value = "Provincia de Limarí"
$ws.Range("S181").Value = 3643
$ws.Range("T181").Value = 1
